$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.891.62"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "2.817.58"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  -0.01%  "
$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.45"
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = "  -0.21%  "
$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.13"
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = "  +4.75%  "
$origStyle_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.561"
$ws.Range("D7").Style = $origStyle_D7
$ws.Range("E7").Value = "  +2.24%  "
$origStyle_D8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = $origStyle_D8
$ws.Range("E8").Value = "  -0.02%  "
$origStyle_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("D9").Style = $origStyle_D9
$ws.Range("E9").Value = "  +6.45%  "
$origStyle_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.40"
$ws.Range("D10").Style = $origStyle_D10
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("E11").Value = "  -0.79%  "
$origStyle_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0843"
$ws.Range("D12").Style = $origStyle_D12
$ws.Range("E12").Value = "  +1.07%  "
$origStyle_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.91"
$ws.Range("D13").Style = $origStyle_D13
$ws.Range("E13").Value = "  -0.11%  "
$origStyle_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.82"
$ws.Range("D14").Style = $origStyle_D14
$ws.Range("E14").Value = "  +4.04%  "
$ws.Range("D15").Value = "3.256.25"
$ws.Range("E15").Value = "  +1.86%  "
$origStyle_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.967"
$ws.Range("D16").Style = $origStyle_D16
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").Value = "2.810.69"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "51.925.07"
$ws.Range("E18").Value = "  +1.59%  "
$origStyle_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.36"
$ws.Range("D19").Style = $origStyle_D19
$ws.Range("E19").Value = "  +9.16%  "
$origStyle_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.65"
$ws.Range("D20").Style = $origStyle_D20
$ws.Range("E20").Value = "  -0.47%  "
$origStyle_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.63"
$ws.Range("D21").Style = $origStyle_D21
$ws.Range("E21").Value = "  +4.25%  "
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  +1.58%  "
$origStyle_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.59"
$ws.Range("D23").Style = $origStyle_D23
$ws.Range("E23").Value = "  +1.33%  "
$origStyle_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.57"
$ws.Range("D24").Style = $origStyle_D24
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$origStyle_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.26"
$ws.Range("D26").Style = $origStyle_D26
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$origStyle_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = $origStyle_D27
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +0.59%  "
$origStyle_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.93"
$ws.Range("D29").Style = $origStyle_D29
$ws.Range("E29").Value = "  +12.64%  "
$origStyle_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.53"
$ws.Range("D30").Style = $origStyle_D30
$ws.Range("E30").Value = "  +3.60%  "
$origStyle_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.28"
$ws.Range("D31").Style = $origStyle_D31
$ws.Range("E31").Value = "  +1.31%  "
$origStyle_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.76"
$ws.Range("D32").Style = $origStyle_D32
$ws.Range("E32").Value = "  +1.88%  "
$origStyle_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.17"
$ws.Range("D33").Style = $origStyle_D33
$ws.Range("E33").Value = "  +2.22%  "
$origStyle_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0907"
$ws.Range("D34").Style = $origStyle_D34
$ws.Range("E34").Value = "  +9.38%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$origStyle_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.72"
$ws.Range("D35").Style = $origStyle_D35
$ws.Range("E35").Value = "  +4.29%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$origStyle_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0453"
$ws.Range("D36").Style = $origStyle_D36
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("E37").Value = "  -0.09%  "
$origStyle_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.95"
$ws.Range("D38").Style = $origStyle_D38
$ws.Range("E38").Value = "  +3.30%  "
$origStyle_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("D39").Style = $origStyle_D39
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("E41").Value = "  +2.28%  "
$origStyle_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.54"
$ws.Range("D42").Style = $origStyle_D42
$ws.Range("E42").Value = "  +1.24%  "
$origStyle_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "121.67"
$ws.Range("D43").Style = $origStyle_D43
$ws.Range("E43").Value = "  +0.88%  "
$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.24"
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = "  +2.19%  "
$origStyle_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.08"
$ws.Range("D45").Style = $origStyle_D45
$ws.Range("E45").Value = "  -0.70%  "
$origStyle_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.53"
$ws.Range("D46").Style = $origStyle_D46
$ws.Range("E46").Value = "  +8.97%  "
$origStyle_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.46"
$ws.Range("D47").Style = $origStyle_D47
$ws.Range("E47").Value = "  +8.93%  "
$ws.Range("D48").Value = "2.139.47"
$ws.Range("E48").Value = "  +2.48%  "
$origStyle_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.03"
$ws.Range("D49").Style = $origStyle_D49
$ws.Range("E49").Value = "  +12.06%  "
$ws.Range("B50").Value = "BEAM"
$ws.Range("C50").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$origStyle_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0324"
$ws.Range("D50").Style = $origStyle_D50
$ws.Range("E50").Value = "  +16.65%  "
$ws.Range("E51").Value = "  +18.27%  "
